$p = $ppt.ActivePresentation

function Update-Title([int]$slideIndex, [string]$middleText, [int]$oldSegLen) {
    $s  = $p.Slides.Item($slideIndex)
    $sh = $s.Shapes.Item(1)
    $tr = $sh.TextFrame.TextRange
    $full = $tr.Text

    # Locate "IA" inside "...com IA e Python" for this title run.
    $marker = "IA e Python"
    $idx0 = $full.IndexOf($marker)   # 0-based
    $start = $idx0 + 1               # Characters() is 1-based

    # Replace the "IA" (+ optional trailing space) segment with the new middle text.
    $seg = $tr.Characters($start, $oldSegLen)
    $seg.Text = $middleText
}

# Slide 1 (title slide): "Automação de Relatórios com " + "ChatGPT" + " e Python"
Update-Title 1 "ChatGPT" 2

# Slide 27 (closing slide, duplicate of the title slide):
# "Automação de Relatórios com " + "ChatGPT " + "e Python"
Update-Title 27 "ChatGPT " 3
